$d = $word.ActiveDocument
$d.TrackRevisions = $true
$p = $d.Paragraphs.Item(1)
$r = $p.Range
$r2 = $d.Range($r.End, $r.End)
$r2.InsertAfter(" nee")
$d.TrackRevisions = $false
$d.Revisions.AcceptAll()
